$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 6.670167666666667
$ws.Range("H2").Value2 = 20.010503
$ws.Range("I2").Value2 = 0.0423069620011633
$ws.Range("J2").Value2 = 0.0423069620011633
$ws.Range("M2").Value2 = 3.867218333333334
$ws.Range("N2").Value2 = 11.601655
$ws.Range("O2").Value2 = 0.1566152977872902
$ws.Range("P2").Value2 = 0.1566152977872902
$ws.Range("Q2").Value2 = 25.79499468694056
$ws.Range("R2").Value2 = 232.154952182465
$ws.Range("S2").Value2 = 0.006625917452287762
$ws.Range("T2").Value2 = 0.006625917452287762

# Row 3
$ws.Range("G3").Value2 = 6.670167666666667
$ws.Range("H3").Value2 = 20.010503
$ws.Range("I3").Value2 = 0.0423069620011633
$ws.Range("J3").Value2 = 0.0423069620011633
$ws.Range("N3").Value2 = 33.813685
$ws.Range("O3").Value2 = 0.4564642152831324
$ws.Range("P3").Value2 = 0.4564642152831324
$ws.Range("Q3").Value2 = 75.18098279261723
$ws.Range("R3").Value2 = 676.6288451335549
$ws.Range("S3").Value2 = 0.0193116142108743
$ws.Range("T3").Value2 = 0.0193116142108743

# Row 4
$ws.Range("G4").Value2 = 6.670167666666667
$ws.Range("H4").Value2 = 20.010503
$ws.Range("I4").Value2 = 0.0423069620011633
$ws.Range("J4").Value2 = 0.0423069620011633
$ws.Range("M4").Value2 = 5.654344666666667
$ws.Range("N4").Value2 = 16.963034
$ws.Range("O4").Value2 = 0.2289906587711778
$ws.Range("P4").Value2 = 0.2289906587711778
$ws.Range("Q4").Value2 = 37.71542697178911
$ws.Range("R4").Value2 = 339.438842746102
$ws.Range("S4").Value2 = 0.00968789909925357
$ws.Range("T4").Value2 = 0.00968789909925357

# Row 5
$ws.Range("G5").Value2 = 6.670167666666667
$ws.Range("H5").Value2 = 20.010503
$ws.Range("I5").Value2 = 0.0423069620011633
$ws.Range("J5").Value2 = 0.0423069620011633
$ws.Range("M5").Value2 = 0.819389
$ws.Range("N5").Value2 = 2.458167
$ws.Range("O5").Value2 = 0.03318376186120772
$ws.Range("P5").Value2 = 0.03318376186120772
$ws.Range("Q5").Value2 = 5.465462014222334
$ws.Range("R5").Value2 = 49.189158128001
$ws.Range("S5").Value2 = 0.001403904152117767
$ws.Range("T5").Value2 = 0.001403904152117767

# Row 6
$ws.Range("G6").Value2 = 6.670167666666667
$ws.Range("H6").Value2 = 20.010503
$ws.Range("I6").Value2 = 0.0423069620011633
$ws.Range("J6").Value2 = 0.0423069620011633
$ws.Range("M6").Value2 = 3.080288333333333
$ws.Range("N6").Value2 = 9.240864999999999
$ws.Range("O6").Value2 = 0.1247460662971919
$ws.Range("P6").Value2 = 0.1247460662971919
$ws.Range("Q6").Value2 = 20.54603964501056
$ws.Range("R6").Value2 = 184.914356805095
$ws.Range("S6").Value2 = 0.005277627086629894
$ws.Range("T6").Value2 = 0.005277627086629894

# Row 7
$ws.Range("I7").Value2 = 0.9513278459982415
$ws.Range("J7").Value2 = 0.9513278459982416
$ws.Range("M7").Value2 = 3.867218333333334
$ws.Range("N7").Value2 = 11.601655
$ws.Range("O7").Value2 = 0.1566152977872902
$ws.Range("P7").Value2 = 0.1566152977872902
$ws.Range("Q7").Value2 = 580.0344806698363
$ws.Range("R7").Value2 = 5220.310326028526
$ws.Range("S7").Value2 = 0.148992493894356
$ws.Range("T7").Value2 = 0.148992493894356

# Row 8
$ws.Range("I8").Value2 = 0.9513278459982415
$ws.Range("J8").Value2 = 0.9513278459982416
$ws.Range("N8").Value2 = 33.813685
$ws.Range("O8").Value2 = 0.4564642152831324
$ws.Range("P8").Value2 = 0.4564642152831324
$ws.Range("Q8").Value2 = 1690.54356628502
$ws.Range("R8").Value2 = 15214.89209656518
$ws.Range("S8").Value2 = 0.4342471187005799
$ws.Range("T8").Value2 = 0.4342471187005799

# Row 9
$ws.Range("I9").Value2 = 0.9513278459982415
$ws.Range("J9").Value2 = 0.9513278459982416
$ws.Range("M9").Value2 = 5.654344666666667
$ws.Range("N9").Value2 = 16.963034
$ws.Range("O9").Value2 = 0.2289906587711778
$ws.Range("P9").Value2 = 0.2289906587711778
$ws.Range("Q9").Value2 = 848.0811243546524
$ws.Range("R9").Value2 = 7632.730119191871
$ws.Range("S9").Value2 = 0.2178451901625029
$ws.Range("T9").Value2 = 0.2178451901625029

# Row 10
$ws.Range("I10").Value2 = 0.9513278459982415
$ws.Range("J10").Value2 = 0.9513278459982416
$ws.Range("M10").Value2 = 0.819389
$ws.Range("N10").Value2 = 2.458167
$ws.Range("O10").Value2 = 0.03318376186120772
$ws.Range("P10").Value2 = 0.03318376186120772
$ws.Range("Q10").Value2 = 122.8981226596317
$ws.Range("R10").Value2 = 1106.083103936685
$ws.Range("S10").Value2 = 0.03156863669354134
$ws.Range("T10").Value2 = 0.03156863669354135

# Row 11
$ws.Range("I11").Value2 = 0.9513278459982415
$ws.Range("J11").Value2 = 0.9513278459982416
$ws.Range("M11").Value2 = 3.080288333333333
$ws.Range("N11").Value2 = 9.240864999999999
$ws.Range("O11").Value2 = 0.1247460662971919
$ws.Range("P11").Value2 = 0.1247460662971919
$ws.Range("Q11").Value2 = 462.0048028677862
$ws.Range("R11").Value2 = 4158.043225810075
$ws.Range("S11").Value2 = 0.1186744065472614
$ws.Range("T11").Value2 = 0.1186744065472614

# Row 12
$ws.Range("G12").Value2 = 1.003544
$ws.Range("H12").Value2 = 3.010632
$ws.Range("I12").Value2 = 0.0063651920005952
$ws.Range("J12").Value2 = 0.006365192000595201
$ws.Range("M12").Value2 = 3.867218333333334
$ws.Range("N12").Value2 = 11.601655
$ws.Range("O12").Value2 = 0.1566152977872902
$ws.Range("P12").Value2 = 0.1566152977872902
$ws.Range("Q12").Value2 = 3.880923755106667
$ws.Range("R12").Value2 = 34.92831379596001
$ws.Range("S12").Value2 = 0.000996886440646495
$ws.Range("T12").Value2 = 0.0009968864406464952

# Row 13
$ws.Range("G13").Value2 = 1.003544
$ws.Range("H13").Value2 = 3.010632
$ws.Range("I13").Value2 = 0.0063651920005952
$ws.Range("J13").Value2 = 0.006365192000595201
$ws.Range("N13").Value2 = 33.813685
$ws.Range("O13").Value2 = 0.4564642152831324
$ws.Range("P13").Value2 = 0.4564642152831324
$ws.Range("Q13").Value2 = 11.31117356654667
$ws.Range("R13").Value2 = 101.80056209892
$ws.Range("S13").Value2 = 0.002905482371678159
$ws.Range("T13").Value2 = 0.00290548237167816

# Row 14
$ws.Range("G14").Value2 = 1.003544
$ws.Range("H14").Value2 = 3.010632
$ws.Range("I14").Value2 = 0.0063651920005952
$ws.Range("J14").Value2 = 0.006365192000595201
$ws.Range("M14").Value2 = 5.654344666666667
$ws.Range("N14").Value2 = 16.963034
$ws.Range("O14").Value2 = 0.2289906587711778
$ws.Range("P14").Value2 = 0.2289906587711778
$ws.Range("Q14").Value2 = 5.674383664165333
$ws.Range("R14").Value2 = 51.06945297748801
$ws.Range("S14").Value2 = 0.001457569509421326
$ws.Range("T14").Value2 = 0.001457569509421326

# Row 15
$ws.Range("G15").Value2 = 1.003544
$ws.Range("H15").Value2 = 3.010632
$ws.Range("I15").Value2 = 0.0063651920005952
$ws.Range("J15").Value2 = 0.006365192000595201
$ws.Range("M15").Value2 = 0.819389
$ws.Range("N15").Value2 = 2.458167
$ws.Range("O15").Value2 = 0.03318376186120772
$ws.Range("P15").Value2 = 0.03318376186120772
$ws.Range("Q15").Value2 = 0.8222929146160001
$ws.Range("R15").Value2 = 7.400636231544
$ws.Range("S15").Value2 = 0.0002112210155486155
$ws.Range("T15").Value2 = 0.0002112210155486155

# Row 16
$ws.Range("G16").Value2 = 1.003544
$ws.Range("H16").Value2 = 3.010632
$ws.Range("I16").Value2 = 0.0063651920005952
$ws.Range("J16").Value2 = 0.006365192000595201
$ws.Range("M16").Value2 = 3.080288333333333
$ws.Range("N16").Value2 = 9.240864999999999
$ws.Range("O16").Value2 = 0.1247460662971919
$ws.Range("P16").Value2 = 0.1247460662971919
$ws.Range("Q16").Value2 = 3.091204875186667
$ws.Range("R16").Value2 = 27.82084387668
$ws.Range("S16").Value2 = 0.0007940326633006042
$ws.Range("T16").Value2 = 0.0007940326633006043
